$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44389   # D2
$ws.Cells.Item(2, 10).Value = 81   # J2
$ws.Cells.Item(2, 11).Value = 2800   # K2
$ws.Cells.Item(2, 12).Value = 3000   # L2
$ws.Cells.Item(2, 13).Value = 2889   # M2
$ws.Cells.Item(2, 16).Value = 963   # P2

$ws.Cells.Item(3, 4).Value = 44222   # D3
$ws.Cells.Item(3, 10).Value = 45   # J3
$ws.Cells.Item(3, 11).Value = 3000   # K3
$ws.Cells.Item(3, 12).Value = 3000   # L3
$ws.Cells.Item(3, 13).Value = 3000   # M3
$ws.Cells.Item(3, 16).Value = 1000   # P3

$ws.Cells.Item(4, 4).Value = 44537   # D4
$ws.Cells.Item(4, 10).Value = 88   # J4
$ws.Cells.Item(4, 11).Value = 2000   # K4
$ws.Cells.Item(4, 12).Value = 2200   # L4
$ws.Cells.Item(4, 13).Value = 2091   # M4
$ws.Cells.Item(4, 16).Value = 697   # P4

$ws.Cells.Item(5, 4).Value = 44627   # D5
$ws.Cells.Item(5, 10).Value = 78   # J5
$ws.Cells.Item(5, 11).Value = 3500   # K5
$ws.Cells.Item(5, 12).Value = 3500   # L5
$ws.Cells.Item(5, 13).Value = 3500   # M5
$ws.Cells.Item(5, 16).Value = 1167   # P5

$ws.Cells.Item(6, 4).Value = 44559   # D6
$ws.Cells.Item(6, 10).Value = 68   # J6
$ws.Cells.Item(6, 11).Value = 2000   # K6
$ws.Cells.Item(6, 12).Value = 2000   # L6
$ws.Cells.Item(6, 13).Value = 2000   # M6
$ws.Cells.Item(6, 16).Value = 667   # P6

$ws.Cells.Item(7, 4).Value = 44225   # D7
$ws.Cells.Item(7, 10).Value = 56   # J7
$ws.Cells.Item(7, 11).Value = 3000   # K7
$ws.Cells.Item(7, 12).Value = 3000   # L7
$ws.Cells.Item(7, 13).Value = 3000   # M7
$ws.Cells.Item(7, 16).Value = 1000   # P7

$ws.Cells.Item(8, 4).Value = 44193   # D8
$ws.Cells.Item(8, 10).Value = 70   # J8
$ws.Cells.Item(8, 11).Value = 3000   # K8
$ws.Cells.Item(8, 12).Value = 3000   # L8
$ws.Cells.Item(8, 13).Value = 3000   # M8
$ws.Cells.Item(8, 16).Value = 1000   # P8

$ws.Cells.Item(9, 4).Value = 44223   # D9
$ws.Cells.Item(9, 10).Value = 80   # J9
$ws.Cells.Item(9, 11).Value = 2500   # K9
$ws.Cells.Item(9, 12).Value = 3000   # L9
$ws.Cells.Item(9, 13).Value = 2781   # M9
$ws.Cells.Item(9, 16).Value = 927   # P9

$ws.Cells.Item(10, 4).Value = 44804   # D10
$ws.Cells.Item(10, 10).Value = 85   # J10
$ws.Cells.Item(10, 11).Value = 3000   # K10
$ws.Cells.Item(10, 12).Value = 3000   # L10
$ws.Cells.Item(10, 13).Value = 3000   # M10
$ws.Cells.Item(10, 16).Value = 1000   # P10

$ws.Cells.Item(11, 4).Value = 44292   # D11
$ws.Cells.Item(11, 10).Value = 40   # J11
$ws.Cells.Item(11, 11).Value = 3000   # K11
$ws.Cells.Item(11, 12).Value = 3000   # L11
$ws.Cells.Item(11, 13).Value = 3000   # M11
$ws.Cells.Item(11, 16).Value = 1000   # P11

$ws.Cells.Item(12, 4).Value = 44165   # D12
$ws.Cells.Item(12, 10).Value = 68   # J12
$ws.Cells.Item(12, 11).Value = 3000   # K12
$ws.Cells.Item(12, 12).Value = 3000   # L12
$ws.Cells.Item(12, 13).Value = 3000   # M12
$ws.Cells.Item(12, 16).Value = 1000   # P12

$ws.Cells.Item(13, 4).Value = 44187   # D13
$ws.Cells.Item(13, 10).Value = 65   # J13
$ws.Cells.Item(13, 11).Value = 3000   # K13
$ws.Cells.Item(13, 12).Value = 3000   # L13
$ws.Cells.Item(13, 13).Value = 3000   # M13
$ws.Cells.Item(13, 16).Value = 1000   # P13

$ws.Cells.Item(14, 4).Value = 44221   # D14
$ws.Cells.Item(14, 10).Value = 50   # J14
$ws.Cells.Item(14, 11).Value = 2500   # K14
$ws.Cells.Item(14, 12).Value = 2500   # L14
$ws.Cells.Item(14, 13).Value = 2500   # M14
$ws.Cells.Item(14, 16).Value = 833   # P14

$ws.Cells.Item(15, 4).Value = 44536   # D15
$ws.Cells.Item(15, 10).Value = 125   # J15
$ws.Cells.Item(15, 11).Value = 2200   # K15
$ws.Cells.Item(15, 12).Value = 2200   # L15
$ws.Cells.Item(15, 13).Value = 2200   # M15
$ws.Cells.Item(15, 16).Value = 733   # P15

$ws.Cells.Item(16, 4).Value = 44756   # D16
$ws.Cells.Item(16, 10).Value = 104   # J16
$ws.Cells.Item(16, 11).Value = 2800   # K16
$ws.Cells.Item(16, 12).Value = 3000   # L16
$ws.Cells.Item(16, 13).Value = 2904   # M16
$ws.Cells.Item(16, 16).Value = 968   # P16

$ws.Cells.Item(17, 4).Value = 44845   # D17
$ws.Cells.Item(17, 10).Value = 80   # J17
$ws.Cells.Item(17, 11).Value = 2500   # K17
$ws.Cells.Item(17, 12).Value = 2500   # L17
$ws.Cells.Item(17, 13).Value = 2500   # M17
$ws.Cells.Item(17, 16).Value = 833   # P17

$ws.Cells.Item(18, 4).Value = 44166   # D18
$ws.Cells.Item(18, 10).Value = 45   # J18
$ws.Cells.Item(18, 11).Value = 2500   # K18
$ws.Cells.Item(18, 12).Value = 2500   # L18
$ws.Cells.Item(18, 13).Value = 2500   # M18
$ws.Cells.Item(18, 16).Value = 833   # P18

$ws.Cells.Item(19, 4).Value = 44965   # D19
$ws.Cells.Item(19, 10).Value = 87   # J19
$ws.Cells.Item(19, 11).Value = 3000   # K19
$ws.Cells.Item(19, 12).Value = 3000   # L19
$ws.Cells.Item(19, 13).Value = 3000   # M19
$ws.Cells.Item(19, 16).Value = 1000   # P19

$ws.Cells.Item(20, 4).Value = 45118   # D20
$ws.Cells.Item(20, 10).Value = 67   # J20
$ws.Cells.Item(20, 11).Value = 3000   # K20
$ws.Cells.Item(20, 12).Value = 3000   # L20
$ws.Cells.Item(20, 13).Value = 3000   # M20
$ws.Cells.Item(20, 16).Value = 1000   # P20

$ws.Cells.Item(21, 4).Value = 44935   # D21
$ws.Cells.Item(21, 10).Value = 78   # J21
$ws.Cells.Item(21, 11).Value = 3000   # K21
$ws.Cells.Item(21, 12).Value = 3000   # L21
$ws.Cells.Item(21, 13).Value = 3000   # M21
$ws.Cells.Item(21, 16).Value = 1000   # P21

$ws.Cells.Item(22, 4).Value = 44669   # D22
$ws.Cells.Item(22, 10).Value = 92   # J22
$ws.Cells.Item(22, 11).Value = 2500   # K22
$ws.Cells.Item(22, 12).Value = 3000   # L22
$ws.Cells.Item(22, 13).Value = 2755   # M22
$ws.Cells.Item(22, 16).Value = 918   # P22

$ws.Cells.Item(23, 4).Value = 44574   # D23
$ws.Cells.Item(23, 10).Value = 50   # J23
$ws.Cells.Item(23, 11).Value = 3000   # K23
$ws.Cells.Item(23, 12).Value = 3000   # L23
$ws.Cells.Item(23, 13).Value = 3000   # M23
$ws.Cells.Item(23, 16).Value = 1000   # P23

$ws.Cells.Item(24, 4).Value = 44243   # D24
$ws.Cells.Item(24, 10).Value = 45   # J24
$ws.Cells.Item(24, 11).Value = 3000   # K24
$ws.Cells.Item(24, 12).Value = 3000   # L24
$ws.Cells.Item(24, 13).Value = 3000   # M24
$ws.Cells.Item(24, 16).Value = 1000   # P24

$ws.Cells.Item(25, 4).Value = 44179   # D25
$ws.Cells.Item(25, 10).Value = 78   # J25
$ws.Cells.Item(25, 11).Value = 3000   # K25
$ws.Cells.Item(25, 12).Value = 3000   # L25
$ws.Cells.Item(25, 13).Value = 3000   # M25
$ws.Cells.Item(25, 16).Value = 1000   # P25

$ws.Cells.Item(26, 4).Value = 44224   # D26
$ws.Cells.Item(26, 10).Value = 67   # J26
$ws.Cells.Item(26, 11).Value = 3000   # K26
$ws.Cells.Item(26, 12).Value = 3000   # L26
$ws.Cells.Item(26, 13).Value = 3000   # M26
$ws.Cells.Item(26, 16).Value = 1000   # P26

$ws.Cells.Item(27, 4).Value = 45092   # D27
$ws.Cells.Item(27, 10).Value = 90   # J27
$ws.Cells.Item(27, 11).Value = 3000   # K27
$ws.Cells.Item(27, 12).Value = 3500   # L27
$ws.Cells.Item(27, 13).Value = 3278   # M27
$ws.Cells.Item(27, 16).Value = 1093   # P27

$ws.Cells.Item(28, 4).Value = 44937   # D28
$ws.Cells.Item(28, 10).Value = 68   # J28
$ws.Cells.Item(28, 11).Value = 3500   # K28
$ws.Cells.Item(28, 12).Value = 3500   # L28
$ws.Cells.Item(28, 13).Value = 3500   # M28
$ws.Cells.Item(28, 16).Value = 1167   # P28

$ws.Cells.Item(29, 4).Value = 44992   # D29
$ws.Cells.Item(29, 10).Value = 45   # J29
$ws.Cells.Item(29, 11).Value = 4000   # K29
$ws.Cells.Item(29, 12).Value = 4000   # L29
$ws.Cells.Item(29, 13).Value = 4000   # M29
$ws.Cells.Item(29, 16).Value = 1333   # P29

$ws.Cells.Item(30, 4).Value = 44967   # D30
$ws.Cells.Item(30, 10).Value = 110   # J30
$ws.Cells.Item(30, 11).Value = 3000   # K30
$ws.Cells.Item(30, 12).Value = 3300   # L30
$ws.Cells.Item(30, 13).Value = 3136   # M30
$ws.Cells.Item(30, 16).Value = 1045   # P30

$ws.Cells.Item(31, 4).Value = 44340   # D31
$ws.Cells.Item(31, 10).Value = 54   # J31
$ws.Cells.Item(31, 11).Value = 3000   # K31
$ws.Cells.Item(31, 12).Value = 3000   # L31
$ws.Cells.Item(31, 13).Value = 3000   # M31
$ws.Cells.Item(31, 16).Value = 1000   # P31

$ws.Cells.Item(32, 4).Value = 44291   # D32
$ws.Cells.Item(32, 10).Value = 45   # J32
$ws.Cells.Item(32, 11).Value = 3000   # K32
$ws.Cells.Item(32, 12).Value = 3000   # L32
$ws.Cells.Item(32, 13).Value = 3000   # M32
$ws.Cells.Item(32, 16).Value = 1000   # P32

$ws.Cells.Item(33, 4).Value = 45117   # D33
$ws.Cells.Item(33, 10).Value = 56   # J33
$ws.Cells.Item(33, 11).Value = 3000   # K33
$ws.Cells.Item(33, 12).Value = 3000   # L33
$ws.Cells.Item(33, 13).Value = 3000   # M33
$ws.Cells.Item(33, 16).Value = 1000   # P33

$ws.Cells.Item(34, 4).Value = 44557   # D34
$ws.Cells.Item(34, 10).Value = 104   # J34
$ws.Cells.Item(34, 11).Value = 2000   # K34
$ws.Cells.Item(34, 12).Value = 2500   # L34
$ws.Cells.Item(34, 13).Value = 2260   # M34
$ws.Cells.Item(34, 16).Value = 753   # P34

$ws.Cells.Item(35, 4).Value = 44242   # D35
$ws.Cells.Item(35, 10).Value = 95   # J35
$ws.Cells.Item(35, 11).Value = 2500   # K35
$ws.Cells.Item(35, 12).Value = 3000   # L35
$ws.Cells.Item(35, 13).Value = 2737   # M35
$ws.Cells.Item(35, 16).Value = 912   # P35

$ws.Cells.Item(36, 4).Value = 44390   # D36
$ws.Cells.Item(36, 10).Value = 50   # J36
$ws.Cells.Item(36, 11).Value = 3000   # K36
$ws.Cells.Item(36, 12).Value = 3000   # L36
$ws.Cells.Item(36, 13).Value = 3000   # M36
$ws.Cells.Item(36, 16).Value = 1000   # P36

$ws.Cells.Item(37, 4).Value = 44260   # D37
$ws.Cells.Item(37, 10).Value = 60   # J37
$ws.Cells.Item(37, 11).Value = 3500   # K37
$ws.Cells.Item(37, 12).Value = 3500   # L37
$ws.Cells.Item(37, 13).Value = 3500   # M37
$ws.Cells.Item(37, 16).Value = 1167   # P37
